# Generate Report for Handoff
# Re-points the handoff report at the freshly generated package
# (new source-doc UUID + new handoff-package hash) and refreshes the
# handoff timestamps for both target locales.

$wb = $excel.ActiveWorkbook

$oldMd = "c9480876-a14a-479c-9876-6eaadb09dd54.md"
$newMd = "a6c92608-7ea3-448a-8954-098c794ff927.md"

$newZh = "a6c92608-7ea3-448a-8954-098c794ff927.d0f5b429c1606cf8def8c33c12c50f7e0042db19.zh-cn.xlf"
$newDe = "a6c92608-7ea3-448a-8954-098c794ff927.d0f5b429c1606cf8def8c33c12c50f7e0042db19.de-de.xlf"

$newZhTime = "2016-03-03 15:32:21"
$newDeTime = "2016-03-03 15:32:34"

$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/341befdb4bb9e575c0bb85c297b021d6525ebd65/e2e/$newMd"
$configAddress = "https://github.com/OpenLocalizationTest/oltest/blob/341befdb4bb9e575c0bb85c297b021d6525ebd65/.localization-config"
$zhAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2db2ad8e84edecb8a5f5eee068fe317605193e05/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$newZh"
$deAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bef9c7ba0fcc6f5346e6a0ab58b44f7783dc7cb0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$newDe"

# --- Overview sheet ------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddress, "", "", $newMd) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configAddress, "", "", ".localization-config") | Out-Null

# --- zh-cn sheet -----------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("C2").Value = $newZh
$wsZh.Range("D2").Value = $newZhTime

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdAddress, "", "", $newMd) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhAddress, "", "", $newZh) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $configAddress, "", "", ".localization-config") | Out-Null

# --- de-de sheet -----------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("C2").Value = $newDe
$wsDe.Range("D2").Value = $newDeTime

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdAddress, "", "", $newMd) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deAddress, "", "", $newDe) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $configAddress, "", "", ".localization-config") | Out-Null
